$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.476.19"
$ws.Range("D3").Value = "1.626.79"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "'212.99"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'0.495"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'0.250"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "1.854.09"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.651.71"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "'63.87"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "26.486.89"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "'214.74"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D25").Value = "'148.90"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").Value = "'2.93"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("D35").Value = "1.219.84"
$ws.Range("E35").Value = "  +4.69%  "
$ws.Range("D36").Value = "'2.39"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "'0.0173"
$ws.Range("E37").Value = "  +4.55%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.795"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").Value = "'0.500"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "1.764.04"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'92.89"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").Value = "'54.74"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'7.64"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.410"
$ws.Range("E51").Value = "  +0.07%  "
